$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values for the 3 remaining rows (area1 separated out, rest combined)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 10

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 3

# Remove the old row 4 entirely (was A4=2, B4=1)
$ws.Range("A4:B4").ClearContents()
$ws.Rows.Item(4).Delete()
